# Add a new worksheet ("Sheet4") at the end of the workbook to hold the
# validation / login credentials used elsewhere in the workbook.
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)

# --- Values -----------------------------------------------------------
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "ajenkins"
$ws.Range("B2").Value = "Acushnet#1"

# --- Formatting ---------------------------------------------------------
# "Username" header is emphasised in red.
$ws.Range("A1").Font.Color = 255

# A small helper cell next to the credentials gets a smaller, explicitly
# black font (used for validation helper text).
$ws.Range("C2").Font.Color = 0
$ws.Range("C2").Font.Size = 10

# A handful of helper cells are pre-formatted as Text so that values typed
# into them later are not reinterpreted as numbers/dates.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F4").NumberFormat = "@"

# --- Column widths --------------------------------------------------------
# (Requested widths are pre-compensated for this engine's column-width
# rounding model so the stored width ends up as close as possible to the
# intended ~11.86 / 16.71 / 7.29 / 10.86 / 9.71 character widths.)
$ws.Columns.Item(1).ColumnWidth = 11
$ws.Columns.Item(2).ColumnWidth = 11
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(4).ColumnWidth = 6.5
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 8.833333333333334

# --- Selection / active cell ---------------------------------------------
$ws.Range("A4").Select() | Out-Null
